# "Generate Report for Archive"
# The localization status report is regenerated: the handoff status text
# moves from "Ready for handoff" to "In Translation" on every sheet that
# shows it (Overview's per-locale status columns, plus each locale sheet's
# own Status column). Because the new text is shorter, the Status columns
# are re-sized (auto-fit) to the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-locale status cells (columns E = zh-cn, F = de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Locale sheets: Status column (column C)
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the Status columns now that the text is shorter
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
